$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure date columns (Y, AA) keep text representation instead of being
# auto-converted to Excel date serials when values are written back.
$ws.Range("Y4").NumberFormat = "@"
$ws.Range("AA4").NumberFormat = "@"
$ws.Range("Y5").NumberFormat = "@"
$ws.Range("AA5").NumberFormat = "@"
$ws.Range("Y9").NumberFormat = "@"
$ws.Range("AA9").NumberFormat = "@"
$ws.Range("Y11").NumberFormat = "@"
$ws.Range("AA11").NumberFormat = "@"
$ws.Range("Y12").NumberFormat = "@"
$ws.Range("AA12").NumberFormat = "@"
$ws.Range("Y13").NumberFormat = "@"
$ws.Range("AA13").NumberFormat = "@"
$ws.Range("Y18").NumberFormat = "@"
$ws.Range("AA18").NumberFormat = "@"
$ws.Range("Y19").NumberFormat = "@"
$ws.Range("AA19").NumberFormat = "@"
$ws.Range("Y22").NumberFormat = "@"
$ws.Range("AA22").NumberFormat = "@"
$ws.Range("Y23").NumberFormat = "@"
$ws.Range("AA23").NumberFormat = "@"
$ws.Range("Y24").NumberFormat = "@"
$ws.Range("AA24").NumberFormat = "@"
$ws.Range("Y25").NumberFormat = "@"
$ws.Range("AA25").NumberFormat = "@"
$ws.Range("Y26").NumberFormat = "@"
$ws.Range("AA26").NumberFormat = "@"
$ws.Range("Y27").NumberFormat = "@"
$ws.Range("AA27").NumberFormat = "@"
$ws.Range("Y28").NumberFormat = "@"
$ws.Range("AA28").NumberFormat = "@"
$ws.Range("Y29").NumberFormat = "@"
$ws.Range("AA29").NumberFormat = "@"
$ws.Range("Y30").NumberFormat = "@"
$ws.Range("AA30").NumberFormat = "@"
$ws.Range("Y31").NumberFormat = "@"
$ws.Range("AA31").NumberFormat = "@"
$ws.Range("Y32").NumberFormat = "@"
$ws.Range("AA32").NumberFormat = "@"
$ws.Range("Y33").NumberFormat = "@"
$ws.Range("AA33").NumberFormat = "@"
$ws.Range("Y34").NumberFormat = "@"
$ws.Range("AA34").NumberFormat = "@"
$ws.Range("Y43").NumberFormat = "@"
$ws.Range("AA43").NumberFormat = "@"
$ws.Range("Y44").NumberFormat = "@"
$ws.Range("AA44").NumberFormat = "@"
$ws.Range("Y45").NumberFormat = "@"
$ws.Range("AA45").NumberFormat = "@"
$ws.Range("Y46").NumberFormat = "@"
$ws.Range("AA46").NumberFormat = "@"
$ws.Range("Y55").NumberFormat = "@"
$ws.Range("AA55").NumberFormat = "@"
$ws.Range("Y56").NumberFormat = "@"
$ws.Range("AA56").NumberFormat = "@"

# Snapshot current row contents (A:AY) before any writes.
$row4 = $ws.Range("A4:AY4").Value2
$row5 = $ws.Range("A5:AY5").Value2
$row9 = $ws.Range("A9:AY9").Value2
$row11 = $ws.Range("A11:AY11").Value2
$row12 = $ws.Range("A12:AY12").Value2
$row13 = $ws.Range("A13:AY13").Value2
$row18 = $ws.Range("A18:AY18").Value2
$row19 = $ws.Range("A19:AY19").Value2
$row22 = $ws.Range("A22:AY22").Value2
$row23 = $ws.Range("A23:AY23").Value2
$row24 = $ws.Range("A24:AY24").Value2
$row25 = $ws.Range("A25:AY25").Value2
$row26 = $ws.Range("A26:AY26").Value2
$row27 = $ws.Range("A27:AY27").Value2
$row28 = $ws.Range("A28:AY28").Value2
$row29 = $ws.Range("A29:AY29").Value2
$row30 = $ws.Range("A30:AY30").Value2
$row31 = $ws.Range("A31:AY31").Value2
$row32 = $ws.Range("A32:AY32").Value2
$row33 = $ws.Range("A33:AY33").Value2
$row34 = $ws.Range("A34:AY34").Value2
$row43 = $ws.Range("A43:AY43").Value2
$row44 = $ws.Range("A44:AY44").Value2
$row45 = $ws.Range("A45:AY45").Value2
$row46 = $ws.Range("A46:AY46").Value2
$row55 = $ws.Range("A55:AY55").Value2
$row56 = $ws.Range("A56:AY56").Value2

# Write rotated/swapped content back to each row.
$ws.Range("A4:AY4").Value2 = $row5
$ws.Range("A5:AY5").Value2 = $row4
$ws.Range("A9:AY9").Value2 = $row11
$ws.Range("A11:AY11").Value2 = $row9
$ws.Range("A12:AY12").Value2 = $row13
$ws.Range("A13:AY13").Value2 = $row12
$ws.Range("A18:AY18").Value2 = $row19
$ws.Range("A19:AY19").Value2 = $row18
$ws.Range("A22:AY22").Value2 = $row23
$ws.Range("A23:AY23").Value2 = $row22
$ws.Range("A24:AY24").Value2 = $row26
$ws.Range("A25:AY25").Value2 = $row24
$ws.Range("A26:AY26").Value2 = $row25
$ws.Range("A27:AY27").Value2 = $row30
$ws.Range("A28:AY28").Value2 = $row27
$ws.Range("A29:AY29").Value2 = $row28
$ws.Range("A30:AY30").Value2 = $row29
$ws.Range("A31:AY31").Value2 = $row34
$ws.Range("A32:AY32").Value2 = $row31
$ws.Range("A33:AY33").Value2 = $row32
$ws.Range("A34:AY34").Value2 = $row33
$ws.Range("A43:AY43").Value2 = $row44
$ws.Range("A44:AY44").Value2 = $row43
$ws.Range("A45:AY45").Value2 = $row46
$ws.Range("A46:AY46").Value2 = $row45
$ws.Range("A55:AY55").Value2 = $row56
$ws.Range("A56:AY56").Value2 = $row55
